$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$url = "https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960"

$rows = @(
    @{ r = 7; a = "2024-09-13 19:53:39"; d = "Starting price monitoring. Current price: `$69.99"; e = "2024-09-13"; f = "19:53:39" },
    @{ r = 8; a = "2024-09-13 19:54:02"; d = "Price remains the same: `$69.99"; e = "2024-09-13"; f = "19:54:02" },
    @{ r = 9; a = "2024-09-13 19:54:24"; d = "Price remains the same: `$69.99"; e = "2024-09-13"; f = "19:54:24" }
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.r, 1).Value = $row.a
    $ws.Cells.Item($row.r, 2).Value = "start_monitoring_price"
    $ws.Cells.Item($row.r, 3).Value = $url
    $ws.Cells.Item($row.r, 4).Value = $row.d
    # Plain "yyyy-mm-dd" text gets auto-parsed into a date serial by Value;
    # force literal text with a quote-prefix, then strip the resulting
    # quotePrefix cell style so the cell stays unstyled like its neighbours.
    $ws.Cells.Item($row.r, 5).Value = "'" + $row.e
    $ws.Cells.Item($row.r, 5).ClearFormats()
    $ws.Cells.Item($row.r, 6).Value = $row.f
}
